# "code optimization and add new feature"
#
# The codes reference sheet is reshaped: several header columns are
# dropped/renamed/reordered (13 cols -> 10 cols) and two sample/test rows
# are added below the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start clean - the old sheet had 13 columns (A:M) and only a header row;
# the new layout only needs 10 columns (A:J) so the leftover columns/shared
# strings must go rather than just be overwritten.
$ws.Cells.Clear()

# New header row (row 1)
$headers = @("SDOH Domain", "Resource", "Resource Element", "Code System", "Code", "Description", "Grouping", "Definition", "Notes", "id")
for ($col = 1; $col -le $headers.Length; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# New sample/test data rows (rows 2-3) - same row shape repeated twice
$sampleRow = @("test", "Resource", "Resource Element", "Code System", "Code", "Description", "Grouping", "Definition", "notes", 12)
for ($row = 2; $row -le 3; $row++) {
    for ($col = 1; $col -le $sampleRow.Length; $col++) {
        $ws.Cells.Item($row, $col).Value = $sampleRow[$col - 1]
    }
}

# Reset selection back to the top-left cell
[void]$ws.Range("A1").Select()
